# Updates cryptos list values (price + 1h volume %) per the source diff.
# Numeric-looking price strings (e.g. "522.60", "7.00") are forced to Text
# via NumberFormat "@" before assignment so Excel COM does not silently
# coerce them to floating-point numbers (which would lose trailing zeros).
# ClearFormats() afterwards drops the temporary "@" format again so the
# cell's style matches the untouched cells (style 0, i.e. no explicit format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.070.11"
$ws.Range("E2").Value = "  -2.81%  "
$ws.Range("D3").Value = "2.656.52"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.60"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.03"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.572"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +8.61%  "
$ws.Range("E10").Value = "  -3.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.336"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("D13").Value = "3.124.43"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").Value = "59.101.26"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.22"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000137"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.648.43"
$ws.Range("E17").Value = "  -8.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "340.61"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.39"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.37"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.36"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "63.66"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("E24").Value = "  -2.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "0.0₃0803"
$ws.Range("E27").Value = "  -3.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.11"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.64"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.82"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "148.97"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.16"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.894"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.880"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.49"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.73"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.61"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.620"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.01"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "275.92"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.04%  "
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0538"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.53"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.92%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.80"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.98%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.031.28"
$ws.Range("E49").Value = "  -4.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.07"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("E51").Value = "  -2.67%  "
